# Add "Recall" response type to the bot (reads/recalls stored info, e.g. a
# previously stored goal) alongside the existing "Store" type.
#
# Touches the two "live" sheets:
#   - User_Initiated_Messages: new "Store"/"Recall" columns + a new
#     Recall row (the bot reading back the stored goal).
#   - Follow_Up_Messages: new "Store"/"Recall" columns (Recall left blank
#     here, a Store column is used instead) + one message text tweak.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# User_Initiated_Messages
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("User_Initiated_Messages")

# Insert a new column at G (old "Input"->"Follow Ups" layout becomes
# Store / Recall / Follow Ups); formatting of the surrounding cells is
# inherited automatically by Excel's column insert.
$ws3.Range("G1:G2").EntireColumn.Insert()

# Header row: F1 "Input" -> "Store", new G1 -> "Recall" (H1 "Follow Ups"
# slides right automatically).
$ws3.Range("F1").Value = "Store"
$ws3.Range("G1").Value = "Recall"

# New row describing the Recall flow: user asks to recall their goal,
# bot responds with the stored value.
$ws3.Range("A3").Value = 1
$ws3.Range("B3").Value = "my goal"
$ws3.Range("C3").Value = "Message;Recall"
$ws3.Range("D3").Value = "You set a goal of [Goals].  How is that going?"
$ws3.Range("G3").Value = "Goals"

$null = $ws3.Range("H3").Select()

# ---------------------------------------------------------------------
# Follow_Up_Messages
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Follow_Up_Messages")

# Same column insert as above (old G "Follow Ups" column slides to H).
$ws4.Range("G1:G9").EntireColumn.Insert()

$ws4.Range("F1").Value = "Store"
$ws4.Range("G1").Value = "Recall"

# Row for "Great, so to achieve that..." now also stores the goal.
$ws4.Range("C4").Value = "Message;Buttons;Store"
$ws4.Range("F4").Value = "Benefits"

# Final row: the bot now stores the goal and gives updated copy.
$ws4.Range("C9").Value = "Message;Store"
$ws4.Range("D9").Value = "That’s great! I'll hold you to that!"
$ws4.Range("F9").Value = "Goals"

$null = $ws4.Range("D9").Select()

# ---------------------------------------------------------------------
# Make User_Initiated_Messages the active tab (was Follow_Up_Messages).
# ---------------------------------------------------------------------
$null = $ws3.Activate()
